# Update dSF column (F) values for a set of rows to reflect repulled data.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("F2").Value = -7
$ws.Range("F4").Value = -1
$ws.Range("F5").Value = 1
$ws.Range("F11").Value = -1
$ws.Range("F13").Value = 1
$ws.Range("F22").Value = -3
$ws.Range("F23").Value = 1
$ws.Range("F24").Value = -5
$ws.Range("F30").Value = -3
$ws.Range("F31").Value = -7
$ws.Range("F32").Value = -6
$ws.Range("F35").Value = 0
$ws.Range("F37").Value = -1
$ws.Range("F38").Value = 1
$ws.Range("F43").Value = 4
$ws.Range("F45").Value = 1
$ws.Range("F46").Value = -3
$ws.Range("F50").Value = -1
$ws.Range("F55").Value = -1
$ws.Range("F56").Value = 1
$ws.Range("F57").Value = 2
$ws.Range("F58").Value = 1
$ws.Range("F60").Value = 0
$ws.Range("F65").Value = -2
$ws.Range("F67").Value = 0
$ws.Range("F81").Value = -1
$ws.Range("F82").Value = 4
$ws.Range("F86").Value = 2
$ws.Range("F88").Value = 0
$ws.Range("F89").Value = -3

$wb.Save()
